$d = $word.ActiveDocument

# --- Change 1: "Stanodavac bira da li šalje obaveštenje ili opomenu"
#     becomes "Stanodavac unosi naslov i tekst obaveštenja ili opomene"
$rng = $d.Content
$rng.Find.Execute("Stanodavac bira da li šalje obaveštenje ili opomenu", $true, $false, $false, $false, $false, $true, 1, $false, "Stanodavac unosi naslov i tekst obaveštenja ili opomene", 2) | Out-Null

# --- Change 2: delete the (now duplicate) paragraph
#     "Stanodavac unosi naslov i tekst obaveštenja ili opomene" that
#     originally followed the paragraph we just replaced.
$rng2 = $d.Content
$rng2.Find.Execute("Stanodavac unosi naslov i tekst obaveštenja ili opomene") | Out-Null
# move past the first match (the one we just created) to find the 2nd (original) occurrence
$afterFirst = $d.Range($rng2.End, $d.Content.End)
$rng3 = $afterFirst
$rng3.Find.Execute("Stanodavac unosi naslov i tekst obaveštenja ili opomene") | Out-Null
# rng3 now covers the text of the duplicate paragraph; expand to include paragraph mark and delete
$pStart = $rng3.Start
$pEnd = $rng3.End
$fullPara = $d.Range($pStart, $pEnd + 1)
$fullPara.Delete() | Out-Null

Write-Output "done change1-2"

# --- Change 3: "Stanodavac pritiska dugme "Potvrdi slanje obaveštenja/opmene""
#     becomes "Stanodavac pritiska dugme "Obaveštenje" ili "Opomena""
$rng4 = $d.Content
$rng4.Find.Execute("Potvrdi slanje obaveštenja/opmene”", $true, $false, $false, $false, $false, $true, 1, $false, "Obaveštenje” ili “Opomena”", 2) | Out-Null

Write-Output "done change3"

# --- Change 4: delete the stray tab-only paragraph that sits right before
#     the "2.2.1.4.a Stanodavac nije kliknuo..." heading
$rngH = $d.Content
$rngH.Find.Execute("Stanodavac nije kliknuo na dugme") | Out-Null
$pH = $rngH.Paragraphs(1)
$prevH = $pH.Previous()
$prevH.Range.Delete() | Out-Null

Write-Output "done change4"

# --- Change 5: "Stanodavac nije kliknuo na dugme "Pošalji          obaveštenje/opomenu""
#     becomes "Stanodavac nije kliknuo na dugme "Obaveštenje" ili "Opomena""
$rng5 = $d.Content
$rng5.Find.Execute("Stanodavac nije kliknuo na dugme") | Out-Null
$p5 = $rng5.Paragraphs(1)
$target5 = $p5.Range
$target5.Find.Execute("“Pošalji          obaveštenje/opomenu”", $true, $false, $false, $false, $false, $true, 1, $false, "“Obaveštenje” ili “Opomena”", 2) | Out-Null

Write-Output "done change5"

# --- Change 6: the paragraph with two tab characters that used to follow the
#     "2.2.1.4.a ..." heading gets merged into that heading paragraph, and a
#     brand-new empty paragraph (carrying the relocated "_GoBack" bookmark)
#     is inserted in its place.
$rng6 = $d.Content
$rng6.Find.Execute("Stanodavac nije kliknuo na dugme") | Out-Null
$p6 = $rng6.Paragraphs(1)
$twoTab = $p6.Next()
Write-Output ("twoTab text: [" + $twoTab.Range.Text + "]")

# 1) Insert a brand new (empty) paragraph mark right at the end of the
#    two-tab paragraph's own content (i.e. just before its existing mark).
#    This new paragraph naturally inherits the two-tab paragraph's (Normal)
#    formatting, matching the target's un-styled paragraph.
$twoTabEnd = $twoTab.Range.End
$insPoint6 = $d.Range($twoTabEnd - 1, $twoTabEnd - 1)
$insPoint6.InsertAfter([char]13)

# 2) Merge the heading paragraph with the (now-preceding) two-tab paragraph
#    by deleting the heading paragraph's own mark.
$rng6b = $d.Content
$rng6b.Find.Execute("Stanodavac nije kliknuo na dugme") | Out-Null
$p6b = $rng6b.Paragraphs(1)
$p6bEnd = $p6b.Range.End
$mark6 = $d.Range($p6bEnd - 1, $p6bEnd)
$mark6.Delete()

# 3) The merge above drops the Heading-3 formatting (it takes on the
#    following/two-tab paragraph's "Normal" style) -- restore it.
$rng6c = $d.Content
$rng6c.Find.Execute("Stanodavac nije kliknuo na dugme") | Out-Null
$p6c = $rng6c.Paragraphs(1)
$p6c.Style = "Heading 3"
Write-Output ("merged heading text: [" + $p6c.Range.Text + "] style=" + $p6c.Style.NameLocal)

$goBackPara = $p6c.Next()
Write-Output ("goBackPara text: [" + $goBackPara.Range.Text + "] len=" + $goBackPara.Range.Text.Length + " style=" + $goBackPara.Style.NameLocal)

# 4) Drop a fresh "_GoBack" bookmark into that new empty paragraph.
$goBackPoint = $d.Range($goBackPara.Range.Start, $goBackPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

Write-Output "done change6"
